$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44959
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 5000
$ws.Range("P2").Value = 5000

# Row 3
$ws.Range("D3").Value = 44498
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = 4000
$ws.Range("P3").Value = 4000

# Row 4
$ws.Range("D4").Value = 44313
$ws.Range("J4").Value = 20

# Row 5
$ws.Range("D5").Value = 44259
$ws.Range("J5").Value = 30

# Row 6
$ws.Range("D6").Value = 44176
$ws.Range("J6").Value = 10

# Row 7
$ws.Range("D7").Value = 44956
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = 5000
$ws.Range("P7").Value = 5000

# Row 8
$ws.Range("D8").Value = 44508
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 4000
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 4000
$ws.Range("P8").Value = 4000

# Row 9
$ws.Range("D9").Value = 44966
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 5000
$ws.Range("P9").Value = 5000

# Row 11
$ws.Range("D11").Value = 44649
$ws.Range("J11").Value = 20
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 5000
$ws.Range("P11").Value = 5000

# Row 12
$ws.Range("D12").Value = 44679
$ws.Range("J12").Value = 50

# Row 13
$ws.Range("D13").Value = 44316
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("P13").Value = 4000

# Row 14
$ws.Range("D14").Value = 44365
$ws.Range("J14").Value = 55

# Row 15
$ws.Range("D15").Value = 44390
$ws.Range("K15").Value = 6000
$ws.Range("L15").Value = 6000
$ws.Range("M15").Value = 6000
$ws.Range("P15").Value = 6000

# Row 17
$ws.Range("D17").Value = 44315
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 4000
$ws.Range("L17").Value = 4000
$ws.Range("M17").Value = 4000
$ws.Range("P17").Value = 4000

# Row 18
$ws.Range("D18").Value = 44680
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 5000
$ws.Range("P18").Value = 5000

# Row 19
$ws.Range("D19").Value = 44280
$ws.Range("J19").Value = 55

# Row 20
$ws.Range("D20").Value = 44781
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = 5000
$ws.Range("P20").Value = 5000

# Row 21
$ws.Range("D21").Value = 44291
$ws.Range("J21").Value = 35
$ws.Range("K21").Value = 4000
$ws.Range("L21").Value = 4000
$ws.Range("M21").Value = 4000
$ws.Range("P21").Value = 4000

# Row 22
$ws.Range("D22").Value = 44509
$ws.Range("J22").Value = 20

# Row 23
$ws.Range("D23").Value = 44497
$ws.Range("J23").Value = 20

# Row 24
$ws.Range("D24").Value = 44749
$ws.Range("J24").Value = 65
$ws.Range("K24").Value = 6000
$ws.Range("L24").Value = 6000
$ws.Range("M24").Value = 6000
$ws.Range("P24").Value = 6000

# Row 25
$ws.Range("D25").Value = 44957
$ws.Range("J25").Value = 20

# Row 26
$ws.Range("D26").Value = 44656
$ws.Range("J26").Value = 85
$ws.Range("K26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = 5000
$ws.Range("P26").Value = 5000

# Row 27
$ws.Range("D27").Value = 44777
$ws.Range("J27").Value = 25

# Row 28
$ws.Range("D28").Value = 44301
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 3000
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = 3000
$ws.Range("P28").Value = 3000
